$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.402.57"
$ws.Range("E2").Value = "  -1.07%  "
$ws.Range("D3").Value = "2.049.06"
$ws.Range("E3").Value = "  -1.56%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.613"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.75%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.15"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.47%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0801"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.74%  "
$ws.Range("E11").Value = "  -2.08%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.74"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.63%  "
$ws.Range("D13").Value = "2.351.86"
$ws.Range("E13").Value = "  -1.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.81"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.759"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.11%  "
$ws.Range("E16").Value = "  -1.24%  "
$ws.Range("D17").Value = "2.047.09"
$ws.Range("E17").Value = "  -2.01%  "
$ws.Range("D18").Value = "37.288.85"
$ws.Range("E18").Value = "  -1.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.09"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.79"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.32%  "
$ws.Range("D21").Value = "0.0₃0833"
$ws.Range("E21").Value = "  -1.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "226.57"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.57%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("E24").Value = "  -0.20%  "
$ws.Range("E25").Value = "  -3.99%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.58"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.89%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "168.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.79%  "
$ws.Range("E28").Value = "  -5.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.37"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.01"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.68%  "
$ws.Range("E31").Value = "  -2.65%  "
$ws.Range("E32").Value = "  -4.28%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0616"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.78%  "
$ws.Range("E34").Value = "  -1.67%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.44"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.35%  "
$ws.Range("E36").Value = "  +0.52%  "
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("E38").Value = "  -4.52%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.38"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.17%  "
$ws.Range("E40").Value = "  -4.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.29"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.38%  "
$ws.Range("D42").Value = "1.496.23"
$ws.Range("E42").Value = "  +3.13%  "
$ws.Range("E43").Value = "  -1.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "96.67"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0941"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.64%  "
$ws.Range("E46").Value = "  +0.79%  "
$ws.Range("E47").Value = "  -4.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.92"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.60%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.16"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.57%  "
$ws.Range("E50").Value = "  -2.37%  "
$ws.Range("D51").Value = "2.239.14"
$ws.Range("E51").Value = "  -1.56%  "
